$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE changes from "001" to "003" (keep as text w/ leading zero)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "003"

# N2: REPORT_DATE changes (stored as plain text, not a real date)
$ws.Range("N2").Value = "2020-03-31 00:00:00"

# Numeric financial figures for row 2
$ws.Range("O2").Value = 62657150.56
$ws.Range("P2").Value = 1568513517.58
$ws.Range("Q2").Value = 1510702102.43
$ws.Range("S2").Value = 1385477183.16
$ws.Range("T2").Value = 1385477183.16
$ws.Range("V2").Value = 13139073.38
$ws.Range("W2").Value = 28309431.23
$ws.Range("X2").Value = 13552097.08
$ws.Range("Y2").Value = 69143225.39
$ws.Range("Z2").Value = 73528909.38
$ws.Range("AA2").Value = 10871758.82
$ws.Range("AG2").Value = 5103082.34
$ws.Range("AS2").Value = 36827688.09
